$wb = $excel.ActiveWorkbook
$dataSheet = $wb.Worksheets.Item("data")

# Refresh the query timestamps on the "data" sheet (F2/F3)
$dataSheet.Range("F2").Value = "2021-10-05 14:22:35.509761"
$dataSheet.Range("F3").Value = "2021-10-05 14:22:35.509771"

# Add the new "metadata" sheet right after "data"
$ws = $wb.Worksheets.Add($null, $dataSheet)
$ws.Name = "metadata"

# Copy the bold/boxed header style from the data sheet onto the metadata
# header row (B1:G1) and the index cell (A2) so the formatting matches.
$dataSheet.Range("B1:F1").Copy($ws.Range("B1:F1"))
$dataSheet.Range("F1").Copy($ws.Range("G1"))
$dataSheet.Range("A2").Copy($ws.Range("A2"))

# Header row
$ws.Range("B1").Value = "data_name"
$ws.Range("C1").Value = "data_id"
$ws.Range("D1").Value = "data_version"
$ws.Range("E1").Value = "data_version_created"
$ws.Range("F1").Value = "panel_query_time"
$ws.Range("G1").Value = "panel_get_request"

# Data row
$ws.Range("B2").Value = "Rhabdoid tumour predisposition"
$ws.Range("C2").Value = 600

# data_version must stay text ("1.6"), not be coerced to a number
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "1.6"
$ws.Range("D2").Style = "Normal"

$ws.Range("E2").Value = "2021-03-08T15:40:37.267588Z"
$ws.Range("F2").Value = "2021-10-05 14:22:35.506116"
$ws.Range("G2").Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/600/?format=json"

Write-Output "done"
